$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 144, shifting existing rows 144:148 down to 145:149.
$ws.Rows(144).Insert()

# Populate the new row 144 with its data (same fixed attributes as the
# surrounding Vega Modelo de Temuco / Papaya rows, with its own measurements).
$ws.Range("A144").Value = 10
$ws.Range("B144").Value = "Vega Modelo de Temuco"
$ws.Range("C144").Value = "La Araucanía"
$ws.Range("D144").Value = 45267
$ws.Range("E144").Value = 9
$ws.Range("F144").Value = "Fruta"
$ws.Range("G144").Value = 100108
$ws.Range("H144").Value = "Tropicales y subtropicales"
$ws.Range("I144").Value = 100108004
$ws.Range("J144").Value = "Papaya"
$ws.Range("K144").Value = "Cultivar IV Región"
$ws.Range("L144").Value = "Primera"
$ws.Range("M144").Value = 25
$ws.Range("N144").Value = 2500
$ws.Range("O144").Value = 2500
$ws.Range("P144").Value = 2500
$ws.Range("Q144").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R144").Value = "Provincia del Elquí"
$ws.Range("S144").Value = 2500
$ws.Range("T144").Value = 1
